# Delete the "2009年" data row (row 2). This shifts the "2010年" row
# (currently row 3) up to row 2, and the sheet's used range shrinks
# from A1:U3 to A1:U2.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(2).Delete()
